# Shift every shape on slide 1 straight down by 219456 EMU (17.28 pt),
# leaving their horizontal position, width and height unchanged.
#
# Note: Shape.Top/Left are stored internally as single-precision floats,
# so a plain "+ deltaPoints" occasionally truncates 1 EMU short when
# converted back on save. A tiny sub-EMU nudge (well under the 1-EMU /
# 12700 point granularity) keeps the float-rounded result on the correct
# side of the EMU boundary without perceptibly moving anything.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$deltaPoints = 219456 / 12700
$nudge = 0.5 / 12700

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $sh.Top = $sh.Top + $deltaPoints + $nudge
}
